$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("V100")

# --- Sheet1: rename entities, insert attribute rows ---

# Row 3 stays row 3: "Row 1" -> "Entity1"
$ws1.Range("A3").Value = "Entity1"

# Insert one attribute row after row 3 (new row 4: Attr1)
$ws1.Rows.Item(4).Insert()
$ws1.Range("A4").Value = "Attr1"
$ws1.Range("A4").Font.Italic = $true

# Old row 4 ("Row 2", B=2) is now row 5: "Row 2" -> "Entity2"
$ws1.Range("A5").Value = "Entity2"

# Insert three attribute rows after row 5 (new rows 6,7,8: Attr1, Attr2, Attr3)
$ws1.Rows.Item(6).Insert()
$ws1.Rows.Item(6).Insert()
$ws1.Rows.Item(6).Insert()
$ws1.Range("A6").Value = "Attr1"
$ws1.Range("A6").Font.Italic = $true
$ws1.Range("A7").Value = "Attr2"
$ws1.Range("A7").Font.Italic = $true
$ws1.Range("A8").Value = "Attr3"
$ws1.Range("A8").Font.Italic = $true

# Old row 5 ("Row 3", B=3) is now row 9: "Row 3" -> "Entity3"
$ws1.Range("A9").Value = "Entity3"

# Insert three attribute rows after row 9 (new rows 10,11,12: Attr1, Attr2, Attr3)
$ws1.Rows.Item(10).Insert()
$ws1.Rows.Item(10).Insert()
$ws1.Rows.Item(10).Insert()
$ws1.Range("A10").Value = "Attr1"
$ws1.Range("A10").Font.Italic = $true
$ws1.Range("A11").Value = "Attr2"
$ws1.Range("A11").Font.Italic = $true
$ws1.Range("A12").Value = "Attr3"
$ws1.Range("A12").Font.Italic = $true

# Row 12's C cell should stay empty (no formula fill)
$ws1.Range("C12").ClearContents()
$ws1.Range("C12").ClearFormats()

# Old total row (SUM) and old V100-reference row have shifted down to rows 13 and 14.
# Delete the V100-reference row (row 14) entirely - no longer needed.
$ws1.Rows.Item(14).Delete()

# Update selection to match new last cell
$ws1.Range("C13").Select() | Out-Null

# --- Sheet V100: clear out the old A1 value ---
$ws2.Range("A1").ClearContents() | Out-Null
